$wb = $excel.ActiveWorkbook

# Rename "Description of Crimes" -> "Description of Crimes (Chicago)"
$ws2 = $wb.Worksheets.Item("Description of Crimes")
$ws2.Name = "Description of Crimes (Chicago)"

# Update its selection to E18 (no longer the tab-selected sheet after we add/activate the new one)
$ws2.Range("E18").Select()

# Add a new worksheet "Data Exploration" at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Data Exploration"

# Narrow column A like the target sheet
$ws3.Columns.Item(1).ColumnWidth = 40.5

# Populate the new sheet's data (match the authoring order so the shared-string
# table indices line up with the target workbook)
$ws3.Range("B1").Value = "6340790 x 22"
$ws3.Range("A4").Value = "Total Dimension (2010 - Current)"
$ws3.Range("A1").Value = "Total Dimensions (2001 - Current): "
$ws3.Range("B4").Value = "2264846 x 24"

# Final selection on the new (now active / tab-selected) sheet
$ws3.Range("E6:E7").Select()
